$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item(1)
$wsZh = $wb.Worksheets.Item(2)
$wsDe = $wb.Worksheets.Item(3)

# ---------------------------------------------------------------------------
# 1. Drop all existing hyperlinks on the three sheets up front. They will be
#    re-created (pointing at the same relationship targets so the same rIds
#    get reused) once the cell content/rows are in their final shape.
# ---------------------------------------------------------------------------
$wsOverview.Range("A1").Hyperlinks.Delete()
$wsZh.Range("A1").Hyperlinks.Delete()
$wsDe.Range("A1").Hyperlinks.Delete()

# ---------------------------------------------------------------------------
# 2. Remove the "Handoff transform failed" row (the ac5acb28-...md entry) on
#    every sheet - row 3 on all three sheets - shifting the
#    ".localization-config" row up into row 3.
# ---------------------------------------------------------------------------
$wsOverview.Rows.Item(3).Delete()
$wsZh.Rows.Item(3).Delete()
$wsDe.Rows.Item(3).Delete()

# ---------------------------------------------------------------------------
# 3. New report was (re)generated under a new source-file UUID, with a new
#    handoff package hash and refreshed handoff timestamps.
# ---------------------------------------------------------------------------
$newMd = "eb09edf2-63a0-4334-9084-477b0fcf42db.md"
$zhXlf = "eb09edf2-63a0-4334-9084-477b0fcf42db.0b9dfe11a05459df4dacf9c9a1ce4369927b9f97.zh-cn.xlf"
$deXlf = "eb09edf2-63a0-4334-9084-477b0fcf42db.0b9dfe11a05459df4dacf9c9a1ce4369927b9f97.de-de.xlf"
$zhTime = "2016-01-15 03:22:13"
$deTime = "2016-01-15 03:22:22"

$wsOverview.Range("A2").Value = $newMd
$wsZh.Range("A2").Value = $newMd
$wsDe.Range("A2").Value = $newMd

$wsZh.Range("C2").Value = $zhXlf
$wsZh.Range("D2").Value = $zhTime

$wsDe.Range("C2").Value = $deXlf
$wsDe.Range("D2").Value = $deTime

# ---------------------------------------------------------------------------
# 4. Re-create the hyperlinks. Addresses match the relationship targets that
#    already exist in each sheet's .rels part so the writer reuses rId2/rId3
#    (Overview) and rId2/rId3/rId4 (zh-cn, de-de) instead of minting new ones.
# ---------------------------------------------------------------------------
$mdTarget = "https://github.com/OpenLocalizationTest/oltest/blob/d6d063d05364ff6764e2a41e7f6b89c92943ce14/e2e/c41e5799-c3fa-4b24-83d4-4aabe91cb9c9.md"
$configTarget = "https://github.com/OpenLocalizationTest/oltest/blob/d6d063d05364ff6764e2a41e7f6b89c92943ce14/.localization-config"
$zhXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/31101f81d9f869a0f84c3a8a33ed528131a00475/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/c41e5799-c3fa-4b24-83d4-4aabe91cb9c9.29a36331b66197022052491d6436265c9cb4c8b9.zh-cn.xlf"
$deXlfTarget = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/775aaf1703f978028f752afd797a36ac9cf2ade3/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/c41e5799-c3fa-4b24-83d4-4aabe91cb9c9.29a36331b66197022052491d6436265c9cb4c8b9.de-de.xlf"

# Overview sheet: A2 (report md) + A3 (.localization-config)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdTarget, "", "", $newMd) | Out-Null
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $configTarget, "", "", ".localization-config") | Out-Null

# zh-cn sheet: A2 (report md) + C2 (xlf handoff package) + A3 (.localization-config)
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdTarget, "", "", $newMd) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), $zhXlfTarget, "", "", $zhXlf) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $configTarget, "", "", ".localization-config") | Out-Null

# de-de sheet: A2 (report md) + C2 (xlf handoff package) + A3 (.localization-config)
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdTarget, "", "", $newMd) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), $deXlfTarget, "", "", $deXlf) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $configTarget, "", "", ".localization-config") | Out-Null
